$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (41 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 831.8
$ws.Range("J12").Value = 1698.25
$ws.Range("L12").Value = 1698.25
$ws.Range("N12").Value = -2038.25
$ws.Range("H76").Value = 16314.8
$ws.Range("I76").Value = 7698.6665
$ws.Range("K76").Value = 7698.6665
$ws.Range("M76").Value = -7383.6665
$ws.Range("H79").Value = 16314.8
$ws.Range("I79").Value = 7698.6665
$ws.Range("K79").Value = 7698.6665
$ws.Range("M79").Value = -6606.6665
$ws.Range("H80").Value = 578.5714
$ws.Range("I80").Value = 562.5
$ws.Range("J80").Value = 600
$ws.Range("K80").Value = 1687.5
$ws.Range("L80").Value = 1800
$ws.Range("M80").Value = -689.5
$ws.Range("N80").Value = -3796
$ws.Range("H83").Value = 578.5714
$ws.Range("I83").Value = 562.5
$ws.Range("J83").Value = 600
$ws.Range("K83").Value = 5062.5
$ws.Range("L83").Value = 5400
$ws.Range("M83").Value = -70.5
$ws.Range("N83").Value = -15384
$ws.Range("H100").Value = 1878.091
$ws.Range("I100").Value = 1436
$ws.Range("K100").Value = 1436
$ws.Range("M100").Value = -895
$ws.Range("H137").Value = 1702
$ws.Range("I137").Value = 1498.9231
$ws.Range("J137").Value = 2230
$ws.Range("K137").Value = 4496.7693
$ws.Range("L137").Value = 6690
$ws.Range("M137").Value = -1946.7693
$ws.Range("N137").Value = -11790
$ws.Range("H138").Value = 2291.06
$ws.Range("J138").Value = 2347.275
$ws.Range("L138").Value = 7041.825000000001
$ws.Range("N138").Value = -17321.825

# --- Sheet: ARM (33 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2958.9546
$ws.Range("I45").Value = 1877.0667
$ws.Range("K45").Value = 1877.0667
$ws.Range("M45").Value = -1500.0667
$ws.Range("H63").Value = 5772.6
$ws.Range("I63").Value = 4100
$ws.Range("J63").Value = 8281.5
$ws.Range("K63").Value = 4100
$ws.Range("L63").Value = 8281.5
$ws.Range("M63").Value = -3414
$ws.Range("N63").Value = -9653.5
$ws.Range("H66").Value = 5772.6
$ws.Range("I66").Value = 4100
$ws.Range("J66").Value = 8281.5
$ws.Range("K66").Value = 20500
$ws.Range("L66").Value = 41407.5
$ws.Range("M66").Value = -17068
$ws.Range("N66").Value = -48271.5
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H122").Value = 3434.4092
$ws.Range("I122").Value = 3238
$ws.Range("J122").Value = 3718.111
$ws.Range("K122").Value = 9714
$ws.Range("L122").Value = 11154.333
$ws.Range("M122").Value = -7264
$ws.Range("N122").Value = -16054.333

# --- Sheet: BSM (14 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1343.95
$ws.Range("I86").Value = 1343.95
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 1343.95
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -220.95
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 1343.95
$ws.Range("I89").Value = 1343.95
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 6719.75
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -1103.75
$ws.Range("N89").Value = 0

# --- Sheet: CRP (28 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3909.625
$ws.Range("J31").Value = 4396.0625
$ws.Range("L31").Value = 4396.0625
$ws.Range("N31").Value = -4986.0625
$ws.Range("H34").Value = 3909.625
$ws.Range("J34").Value = 4396.0625
$ws.Range("L34").Value = 4396.0625
$ws.Range("N34").Value = -4800.0625
$ws.Range("H58").Value = 4141.222
$ws.Range("I58").Value = 2779.7058
$ws.Range("K58").Value = 2779.7058
$ws.Range("M58").Value = -2576.7058
$ws.Range("H103").Value = 25395.2
$ws.Range("I103").Value = 2800
$ws.Range("K103").Value = 2800
$ws.Range("M103").Value = -1628
$ws.Range("H111").Value = 54959
$ws.Range("J111").Value = 54959
$ws.Range("L111").Value = 54959
$ws.Range("N111").Value = -63139
$ws.Range("H112").Value = 79973
$ws.Range("J112").Value = 79973
$ws.Range("L112").Value = 79973
$ws.Range("N112").Value = -82927
$ws.Range("H136").Value = 4141.222
$ws.Range("I136").Value = 2779.7058
$ws.Range("K136").Value = 8339.117400000001
$ws.Range("M136").Value = -5789.117400000001

# --- Sheet: GSM (15 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 18120
$ws.Range("I99").Value = 11493.333
$ws.Range("K99").Value = 11493.333
$ws.Range("M99").Value = -9247.333000000001
$ws.Range("H102").Value = 1838.625
$ws.Range("I102").Value = 1483.2354
$ws.Range("J102").Value = 2701.7144
$ws.Range("K102").Value = 1483.2354
$ws.Range("L102").Value = 2701.7144
$ws.Range("M102").Value = 138.7646
$ws.Range("N102").Value = -5945.7144
$ws.Range("H105").Value = 29475
$ws.Range("J105").Value = 29475
$ws.Range("L105").Value = 29475
$ws.Range("N105").Value = -36463

# --- Sheet: LTW (12 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2675.4
$ws.Range("I40").Value = 2761
$ws.Range("K40").Value = 2761
$ws.Range("M40").Value = -2625
$ws.Range("H101").Value = 73244.75
$ws.Range("J101").Value = 73244.75
$ws.Range("L101").Value = 73244.75
$ws.Range("N101").Value = -79734.75
$ws.Range("H136").Value = 4165.8887
$ws.Range("I136").Value = 3803.9092
$ws.Range("K136").Value = 11411.7276
$ws.Range("M136").Value = -8861.7276

# --- Sheet: WVR (30 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 4998.6665
$ws.Range("I23").Value = 4996
$ws.Range("K23").Value = 4996
$ws.Range("M23").Value = -4767
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("N58").ClearContents()
$ws.Range("H122").Value = 5364.6313
$ws.Range("I122").Value = 5345
$ws.Range("K122").Value = 16035
$ws.Range("M122").Value = -13585
$ws.Range("H126").Value = 5610.25
$ws.Range("I126").Value = 5146.4443
$ws.Range("K126").Value = 15439.3329
$ws.Range("M126").Value = -12969.3329
$ws.Range("H132").Value = 4636.154
$ws.Range("I132").Value = 4240.048
$ws.Range("J132").Value = 6299.8
$ws.Range("K132").Value = 12720.144
$ws.Range("L132").Value = 18899.4
$ws.Range("M132").Value = -10190.144
$ws.Range("N132").Value = -23959.4
$ws.Range("H136").Value = 4077.9355
$ws.Range("I136").Value = 2516.1155
$ws.Range("K136").Value = 7548.3465
$ws.Range("M136").Value = -4998.3465
